# Add "MFG Part#" column (H) to the BOM sheet, plus the new manufacturer
# part-number data gathered during the BOM update, and a new "alternative"
# note for the crystal (Y1) row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("H1").Value = "MFG Part#"

# Row 4 - U1 / ATmega32U4-MU
$ws.Range("H4").Value = "ATMEGA32U4-MU"

# Row 5 - P1 / USB_C_Plug
$ws.Range("H5").Value = "USB4155-03-C "

# Row 6 - Y1 / 16MHz crystal (also gets an alternative note in column G)
$ws.Range("G6").Value = "EB3250YA12-16.000M (alternative)"
$ws.Range("H6").Value = "SA534160F35HDT"

# Row 7 - R5,R4,R2,R1 resistors
$ws.Range("H7").Value = "CRG0603F22R "

# Row 8 - R3 resistor
$ws.Range("H8").Value = "CPF0603B10KE "

# Row 9 - F1 fuse
$ws.Range("H9").Value = "CC12H500mA-TR "

# Row 10 - C7,C5,C2,C1 capacitors
$ws.Range("H10").Value = "GRM033R61A104ME15J "

# Row 11 - C6 capacitor
$ws.Range("H11").Value = "GRM033R60J105MEA2J "

# Row 12 - C4,C3 capacitors
$ws.Range("H12").Value = "GRM0335C1H220JA01J "

# Column widths for the new columns (closest reachable values given Excel's
# pixel-quantized column-width grid; targets are 29.45 and 21.11 chars)
$ws.Columns.Item(7).ColumnWidth = 28.6
$ws.Columns.Item(8).ColumnWidth = 20.4

# Leave selection where the author's cursor ended up after editing
$ws.Range("G13").Select() | Out-Null
